# versão final do relatório AoL2023/2024
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the evaluated scores in column C
$ws.Range("C2").Value = "A"
$ws.Range("C6").Value = "B"

# Move the active selection from C3 to C7
$ws.Range("C7").Select()
